# Scheduled runner update: refresh market-derived Leve profit figures
# across all job sheets (values sourced from latest market-board pull).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 321.2143
$ws.Range("I96").Value = 313.45456
$ws.Range("J96").Value = 349.66666
$ws.Range("K96").Value = 940.36368
$ws.Range("L96").Value = 1048.99998
$ws.Range("M96").Value = 432.63632
$ws.Range("N96").Value = -3794.99998
$ws.Range("H98").Value = 1091.1111
$ws.Range("I98").Value = 1227.7142
$ws.Range("J98").Value = 613
$ws.Range("K98").Value = 1227.7142
$ws.Range("L98").Value = 613
$ws.Range("M98").Value = 270.2858000000001
$ws.Range("N98").Value = -3609
$ws.Range("H101").Value = 1228.4
$ws.Range("I101").Value = 597.6667
$ws.Range("K101").Value = 1793.0001
$ws.Range("M101").Value = -171.0001
$ws.Range("H122").Value = 1091.1111
$ws.Range("I122").Value = 1227.7142
$ws.Range("J122").Value = 613
$ws.Range("K122").Value = 3683.1426
$ws.Range("L122").Value = 1839
$ws.Range("M122").Value = -1233.1426
$ws.Range("N122").Value = -6739
$ws.Range("H123").Value = 69744.44500000001
$ws.Range("J123").Value = 69744.44500000001
$ws.Range("L123").Value = 69744.44500000001
$ws.Range("N123").Value = -79544.44500000001
$ws.Range("H125").Value = 1055.6923
$ws.Range("I125").Value = 539.5
$ws.Range("J125").Value = 1881.6
$ws.Range("K125").Value = 4855.5
$ws.Range("L125").Value = 16934.4
$ws.Range("M125").Value = -2395.5
$ws.Range("N125").Value = -21854.4
$ws.Range("H127").Value = 1366.0555
$ws.Range("J127").Value = 2019.5
$ws.Range("L127").Value = 6058.5
$ws.Range("N127").Value = -15978.5
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 922.6667
$ws.Range("J129").Value = 1150
$ws.Range("L129").Value = 3450
$ws.Range("N129").Value = -13450
$ws.Range("H131").Value = 902.1429000000001
$ws.Range("J131").Value = 1051.4286
$ws.Range("L131").Value = 3154.2858
$ws.Range("N131").Value = -13234.2858
$ws.Range("H132").Value = 1115948.1
$ws.Range("I132").Value = 2322.9707
$ws.Range("J132").Value = 4902273.5
$ws.Range("K132").Value = 6968.9121
$ws.Range("L132").Value = 14706820.5
$ws.Range("M132").Value = -4438.9121
$ws.Range("N132").Value = -14711880.5
$ws.Range("H137").Value = 5002447
$ws.Range("I137").Value = 9092804
$ws.Range("J137").Value = 3121.889
$ws.Range("K137").Value = 27278412
$ws.Range("L137").Value = 9365.667000000001
$ws.Range("M137").Value = -27275862
$ws.Range("N137").Value = -14465.667
$ws.Range("H138").Value = 2565933
$ws.Range("I138").Value = 1444.9269
$ws.Range("J138").Value = 6946933.5
$ws.Range("K138").Value = 4334.780699999999
$ws.Range("L138").Value = 20840800.5
$ws.Range("M138").Value = 805.2193000000007
$ws.Range("N138").Value = -20851080.5
$ws.Range("H141").Value = 2013.8387
$ws.Range("I141").Value = 1997.5518
$ws.Range("J141").Value = 2250
$ws.Range("K141").Value = 5992.6554
$ws.Range("L141").Value = 6750
$ws.Range("M141").Value = -812.6553999999996
$ws.Range("N141").Value = -17110

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 62626104
$ws.Range("I61").Value = 77000744
$ws.Range("J61").Value = 336000
$ws.Range("K61").Value = 77000744
$ws.Range("L61").Value = 336000
$ws.Range("M61").Value = -77000532
$ws.Range("N61").Value = -336424
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H74").Value = 6463029
$ws.Range("I74").Value = 8656404
$ws.Range("K74").Value = 8656404
$ws.Range("M74").Value = -8655530
$ws.Range("H77").Value = 6463029
$ws.Range("I77").Value = 8656404
$ws.Range("K77").Value = 43282020
$ws.Range("M77").Value = -43277652
$ws.Range("H97").Value = 1736776
$ws.Range("I97").Value = 2500636.5
$ws.Range("K97").Value = 2500636.5
$ws.Range("M97").Value = -2500140.5
$ws.Range("H102").Value = 5719937.5
$ws.Range("I102").Value = 6808958.5
$ws.Range("J102").Value = 2577.75
$ws.Range("K102").Value = 6808958.5
$ws.Range("L102").Value = 2577.75
$ws.Range("M102").Value = -6807336.5
$ws.Range("N102").Value = -5821.75
$ws.Range("H136").Value = 62626104
$ws.Range("I136").Value = 77000744
$ws.Range("J136").Value = 336000
$ws.Range("K136").Value = 231002232
$ws.Range("L136").Value = 1008000
$ws.Range("M136").Value = -230999682
$ws.Range("N136").Value = -1013100

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 13233.333
$ws.Range("J21").Value = 13233.333
$ws.Range("L21").Value = 13233.333
$ws.Range("N21").Value = -13705.333
$ws.Range("H94").Value = 420
$ws.Range("I94").Value = 283.52942
$ws.Range("K94").Value = 283.52942
$ws.Range("M94").Value = 167.47058
$ws.Range("H137").Value = 55390
$ws.Range("J137").Value = 55390
$ws.Range("L137").Value = 55390
$ws.Range("N137").Value = -65590

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 1408.5
$ws.Range("I123").Value = 362.75
$ws.Range("J123").Value = 3500
$ws.Range("K123").Value = 1088.25
$ws.Range("L123").Value = 10500
$ws.Range("M123").Value = 1361.75
$ws.Range("N123").Value = -15400
$ws.Range("H130").Value = 3181.25
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1336.75
$ws.Range("I30").Value = 1336.75
$ws.Range("K30").Value = 1336.75
$ws.Range("M30").Value = -1228.75
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H122").Value = 3180.3462
$ws.Range("I122").Value = 2419.25
$ws.Range("J122").Value = 3518.611
$ws.Range("K122").Value = 7257.75
$ws.Range("L122").Value = 10555.833
$ws.Range("M122").Value = -4807.75
$ws.Range("N122").Value = -15455.833
$ws.Range("H132").Value = 38911.785
$ws.Range("I132").Value = 18163.516
$ws.Range("J132").Value = 114988.78
$ws.Range("K132").Value = 54490.548
$ws.Range("L132").Value = 344966.34
$ws.Range("M132").Value = -51960.548
$ws.Range("N132").Value = -350026.34

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2354.4883
$ws.Range("I122").Value = 2023
$ws.Range("J122").Value = 2973.2666
$ws.Range("K122").Value = 6069
$ws.Range("L122").Value = 8919.799800000001
$ws.Range("M122").Value = -3619
$ws.Range("N122").Value = -13819.7998
$ws.Range("H123").Value = 42000
$ws.Range("J123").Value = 42000
$ws.Range("L123").Value = 42000
$ws.Range("N123").Value = -51800
$ws.Range("H124").Value = 55000
$ws.Range("J124").Value = 55000
$ws.Range("L124").Value = 55000
$ws.Range("N124").Value = -64820
$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H126").Value = 1342.2
$ws.Range("I126").Value = 970.3333
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 2910.9999
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -440.9998999999998
$ws.Range("N126").Value = -10640
$ws.Range("H130").Value = 50163
$ws.Range("J130").Value = 50163
$ws.Range("L130").Value = 50163
$ws.Range("N130").Value = -60203
$ws.Range("H131").Value = 44857.5
$ws.Range("J131").Value = 44857.5
$ws.Range("L131").Value = 44857.5
$ws.Range("N131").Value = -54937.5
